$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "88.250.14"
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3").Value = "3.073.91"
$ws.Range("E3").Value = "  -3.80%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "209.65"
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("D6").Value = "622.35"
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("E7").Value = "  -3.12%  "
$ws.Range("D8").Value = "0.822"
$ws.Range("E8").Value = "  +18.39%  "
$ws.Range("D9").Value = "1.00"
$ws.Range("D10").Value = "3.068.13"
$ws.Range("E10").Value = "  -3.88%  "
$ws.Range("D11").Value = "0.593"
$ws.Range("E11").Value = "  +4.51%  "
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").Value = "0.0000237"
$ws.Range("E13").Value = "  -5.01%  "
$ws.Range("D14").Value = "5.28"
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("D15").Value = "88.031.36"
$ws.Range("E15").Value = "  -1.85%  "
$ws.Range("D16").Value = "3.637.56"
$ws.Range("E16").Value = "  -3.85%  "
$ws.Range("D17").Value = "31.48"
$ws.Range("E17").Value = "  -2.97%  "
$ws.Range("D18").Value = "3.084.06"
$ws.Range("E18").Value = "  -3.77%  "
$ws.Range("E19").Value = "  -1.71%  "
$ws.Range("D20").Value = "0.0000211"
$ws.Range("E20").Value = "  -5.48%  "
$ws.Range("D21").Value = "13.05"
$ws.Range("D22").Value = "420.54"
$ws.Range("E22").Value = "  -2.54%  "
$ws.Range("D23").Value = "8.18"
$ws.Range("E23").Value = "  -3.53%  "
$ws.Range("E24").Value = "  -3.68%  "
$ws.Range("D25").Value = "5.44"
$ws.Range("E25").Value = "  +7.85%  "
$ws.Range("D26").Value = "82.05"
$ws.Range("E26").Value = "  +9.02%  "
$ws.Range("D27").Value = "11.35"
$ws.Range("E27").Value = "  -1.27%  "
$ws.Range("D28").Value = "3.236.74"
$ws.Range("E28").Value = "  -3.85%  "
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E30").Value = "  +8.53%  "
$ws.Range("D31").Value = "0.149"
$ws.Range("E31").Value = "  -2.21%  "
$ws.Range("D32").Value = "8.01"
$ws.Range("E32").Value = "  -4.05%  "
$ws.Range("D33").Value = "503.57"
$ws.Range("E33").Value = "  -5.55%  "
$ws.Range("D34").Value = "3.55"
$ws.Range("E34").Value = "  -11.14%  "
$ws.Range("D35").Value = "6.61"
$ws.Range("E35").Value = "  -2.43%  "
$ws.Range("D36").Value = "1.80"
$ws.Range("E36").Value = "  -3.94%  "
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("D38").Value = "22.30"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("D39").Value = "0.130"
$ws.Range("E39").Value = "  +3.31%  "
$ws.Range("D40").Value = "22.21"
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").Value = "0.357"
$ws.Range("E43").Value = "  -2.85%  "
$ws.Range("D44").Value = "147.74"
$ws.Range("E44").Value = "  -1.96%  "
$ws.Range("E45").Value = "  -4.11%  "
$ws.Range("E46").Value = "  +9.57%  "
$ws.Range("D47").Value = "43.50"
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("D48").Value = "0.0669"
$ws.Range("E48").Value = "  +13.54%  "
$ws.Range("D49").Value = "156.45"
$ws.Range("E49").Value = "  -8.14%  "
$ws.Range("D50").Value = "0.701"
$ws.Range("E50").Value = "  -3.48%  "
$ws.Range("E51").Value = "  -4.44%  "
